$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.480892062187195
$ws.Range("B1").Value = 4.160904407501221
$ws.Range("C1").Value = 3.523420810699463
$ws.Range("D1").Value = 1.865121126174927
$ws.Range("E1").Value = 0.6341943740844727
